$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.867.74"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.671.35"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'216.16"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "'0.532"
$ws.Range("E6").Value = "  +5.43%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "'0.255"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("D9").Value = "'0.0619"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").Value = "'20.31"
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("D11").Value = "'0.0892"
$ws.Range("E11").Value = "  +3.86%  "
$ws.Range("D12").Value = "1.907.13"
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "1.659.00"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "'65.67"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "26.893.09"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "'232.95"
$ws.Range("E18").Value = "  -3.46%  "
$ws.Range("D19").Value = "'7.85"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'9.21"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.21"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "'145.70"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'0.117"
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "'0.0497"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "1.466.68"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("D34").Value = "'3.18"
$ws.Range("E34").Value = "  +4.74%  "
$ws.Range("E35").Value = "  +3.24%  "
$ws.Range("D36").Value = "'2.41"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "'0.903"
$ws.Range("E37").Value = "  +4.87%  "
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "'5.93"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("E43").Value = "  +6.74%  "
$ws.Range("D44").Value = "'65.78"
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").Value = "1.815.95"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("D47").Value = "'90.43"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  +13.57%  "
$ws.Range("D50").Value = "'0.100"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("E51").Value = "  +1.40%  "
